# Apply wording updates described in the commit diff.
$p = $ppt.ActivePresentation

# --- Slide 3: "Project Idea" - Content Placeholder 2, paragraph 2 ---
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item("Content Placeholder 2")
$tr3 = $sh3.TextFrame.TextRange
$para3_2 = $tr3.Paragraphs(2)
$para3_2.Runs(1).Text = "An Application which allows users to share their opinions about suggested policies."

# --- Slide 4: "Project Description (MVP)" - Content Placeholder 2, paragraphs 1-3 ---
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item("Content Placeholder 2")
$tr4 = $sh4.TextFrame.TextRange

$para4_1 = $tr4.Paragraphs(1)
$para4_1.Runs(1).Text = "This application will allow users to write new policies and share them with others."

$para4_2 = $tr4.Paragraphs(2)
$para4_2.Runs(1).Text = "This application will allow users to read the policies that have been published and vote to agree or disagree."

$para4_3 = $tr4.Paragraphs(3)
$para4_3.Runs(1).Text = "This application will require that users sign in to publish new policies or vote."
